$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.806.00"
$ws.Range("E2").Value = "  -1.92%  "
$ws.Range("D3").Value = "1.889.08"
$ws.Range("E3").Value = "  -1.91%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'0.7702"
$ws.Range("E5").Value = "  -4.76%  "
$ws.Range("D6").Value = "'244.40"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.3124"
$ws.Range("D9").Value = "'25.30"
$ws.Range("E9").Value = "  -7.46%  "
$ws.Range("D10").Value = "'0.07214"
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("D11").Value = "'0.08095"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "'0.7662"
$ws.Range("E12").Value = "  -4.30%  "
$ws.Range("D13").Value = "'5.536"
$ws.Range("E13").Value = "  +2.03%  "
$ws.Range("D14").Value = "1.883.21"
$ws.Range("E14").Value = "  -2.16%  "
$ws.Range("D15").Value = "'92.24"
$ws.Range("E15").Value = "  -2.54%  "
$ws.Range("D16").Value = "'6.150"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "29.815.22"
$ws.Range("E17").Value = "  -1.86%  "
$ws.Range("D18").Value = "'13.94"
$ws.Range("E18").Value = "  -3.17%  "
$ws.Range("D19").Value = "'243.20"
$ws.Range("E19").Value = "  -4.02%  "
$ws.Range("D20").Value = "'0.000007768"
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").Value = "'8.156"
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("B23").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D23").Value = "2.143.39"
$ws.Range("E23").Value = "  -1.67%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'0.1558"
$ws.Range("E25").Value = "  -4.76%  "
$ws.Range("D26").Value = "'9.425"
$ws.Range("E26").Value = "  -1.64%  "
$ws.Range("D27").Value = "'162.25"
$ws.Range("E27").Value = "  -3.24%  "
$ws.Range("E28").Value = "  -2.34%  "
$ws.Range("D29").Value = "'2.041"
$ws.Range("E29").Value = "  -5.60%  "
$ws.Range("D30").Value = "'1.442"
$ws.Range("E30").Value = "  +4.82%  "
$ws.Range("D31").Value = "'1.549"
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").Value = "'4.459"
$ws.Range("E32").Value = "  +2.35%  "
$ws.Range("D33").Value = "'4.082"
$ws.Range("E33").Value = "  -1.99%  "
$ws.Range("D34").Value = "'0.05499"
$ws.Range("E34").Value = "  -2.42%  "
$ws.Range("D35").Value = "'1.259"
$ws.Range("E35").Value = "  -3.59%  "
$ws.Range("D36").Value = "'0.7470"
$ws.Range("D37").Value = "'1.001"
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("D38").Value = "'2.627"
$ws.Range("E38").Value = "  -3.41%  "
$ws.Range("D39").Value = "'0.01921"
$ws.Range("E39").Value = "  -2.03%  "
$ws.Range("D40").Value = "'2.781"
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("D41").Value = "1.155.55"
$ws.Range("E41").Value = "  +11.71%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'73.58"
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.4414"
$ws.Range("E43").Value = "  -2.24%  "
$ws.Range("D44").Value = "'5.893"
$ws.Range("E44").Value = "  -1.88%  "
$ws.Range("D45").Value = "'0.8507"
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("D46").Value = "'0.9999"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "'103.27"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").Value = "'1.884"
$ws.Range("E48").Value = "  -3.00%  "
$ws.Range("D49").Value = "'9.928"
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("D50").Value = "'7.448"
$ws.Range("E50").Value = "  -3.05%  "
$ws.Range("D51").Value = "'2.999"
$ws.Range("E51").Value = "  +9.74%  "
